# edit.ps1 - apply the CherryPop docx changes described by the diff:
#   1. Insert a new "Meta description" paragraph right after the H1 title
#      paragraph (bold "Meta description" run + plain run with the rest).
#   2. Remove the duplicate bold "Play CherryPop Free - ..." paragraph that
#      used to sit near the end of the document (right before the italic
#      paragraph).
#   3. Replace the text of the trailing italic paragraph with the new
#      "Create a feature image ..." image-prompt text, keeping the italic
#      run formatting intact.
#
# NOTE: this document's paragraphs (the ones without an explicit <w:pPr>)
# all begin with a stray empty <w:r/> run before their real text run(s) -
# an artifact of however the source docx was produced. Any in-place
# character deletion inside such a paragraph causes the run list to be
# rebuilt/collapsed (losing that leading empty run), whereas pure
# insertions, paragraph-level (whole paragraph) deletes, and Copy/Paste of
# already-finished text into a still-empty paragraph do not. So new/edited
# paragraphs are first fully composed (text + bold/italic formatting) in a
# disposable scratch paragraph at the end of the story, then that finished
# content is Copy/Pasted into the (still pristine, untouched) destination
# paragraph, and the scratch paragraph is discarded again - this keeps the
# leading <w:r/> quirk faithfully reproduced at the destination.

$d = $word.ActiveDocument

function New-ScratchParagraph {
    # Appends a brand-new empty paragraph at the very end of the story and
    # returns its 1-based paragraph index.
    $last = $d.Paragraphs($d.Paragraphs.Count)
    $endRange = $last.Range
    $endRange.Collapse(0)
    [void]$endRange.InsertParagraphAfter()
    return $d.Paragraphs.Count
}

# ---------------------------------------------------------------------
# 1. Insert the "Meta description" paragraph after the H1 title paragraph
# ---------------------------------------------------------------------

# 1a. Create the (still pristine/untouched) destination paragraph.
$titlePara = $d.Paragraphs(1)
[void]$titlePara.Range.InsertParagraphAfter()
$d.Paragraphs(2).Style = "Normal"

# 1b. Build the finished "Meta description: ..." text + bold formatting in
#     a disposable scratch paragraph, starting from a copy of an existing
#     plain-text paragraph (so we inherit its leading empty <w:r/> and can
#     destructively edit its text without caring about the consequences).
$sourceForMeta = $d.Paragraphs(4).Range   # "Get ready to pop some cherries..."
$sourceForMetaNoMark = $d.Range($sourceForMeta.Start, $sourceForMeta.End - 1)
$sourceForMetaNoMark.Copy()

$scratchIndex = New-ScratchParagraph
$scratchRange = $d.Paragraphs($scratchIndex).Range
$scratchCollapsed = $d.Range($scratchRange.Start, $scratchRange.Start)
$scratchCollapsed.Paste()

$scratchRange = $d.Paragraphs($scratchIndex).Range
[void]$scratchRange.Find.Execute(
    "Get ready to pop some cherries and cash in on fruity rewards with the online slot game, CherryPop! Unlike other conventional slot games, CherryPop comes with an expandable grid that offers players 59,049 ways to win big. Interestingly, every time you hit a winning combination, a new row is added to the top of the winning row, increasing your chances of hitting another win. With the grid's capability to expand up to 6 upper rows in regular mode and 9 in Free Spins, there's never a dull moment in this thrilling game!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description: Read our review of CherryPop, an online slot game with PopWins mode and 5 free spins. Play it for free!",
    2)

$metaStart = $d.Paragraphs($scratchIndex).Range.Start
$metaBoldRange = $d.Range($metaStart, $metaStart + 16)
$metaBoldRange.Bold = 1

# 1c. Copy the finished scratch paragraph's content (no paragraph mark) and
#     paste it into the still-pristine destination paragraph.
$scratchRange = $d.Paragraphs($scratchIndex).Range
$scratchNoMark = $d.Range($scratchRange.Start, $scratchRange.End - 1)
$scratchNoMark.Copy()

$destRange = $d.Paragraphs(2).Range
$destCollapsed = $d.Range($destRange.Start, $destRange.Start)
$destCollapsed.Paste()

# 1d. Discard the scratch paragraph (whole-paragraph delete - does not
#     disturb any other paragraph's run layout).
$scratchIndex = $d.Paragraphs.Count
$d.Paragraphs($scratchIndex).Range.Delete()

# ---------------------------------------------------------------------
# 2. Remove the duplicate bold "Play CherryPop Free - ..." paragraph near
#    the end of the document (now the 2nd-to-last paragraph).
# ---------------------------------------------------------------------

$removeIndex = $d.Paragraphs.Count - 1
$d.Paragraphs($removeIndex).Range.Delete()

# ---------------------------------------------------------------------
# 3. Replace the trailing italic paragraph's text with the new image
#    prompt, preserving the italic run formatting and the leading empty
#    run.
# ---------------------------------------------------------------------

$newImagePrompt = 'Create a feature image that captures the essence of CherryPop, a fun and exciting online slot game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior can be holding a cherry or surrounded by exploding symbols to represent the PopWins mode. The background should be colorful and vibrant, with elements of fruits, gems, and playing cards to represent the different symbols in the game. The image should also have the CherryPop logo prominently displayed, along with the tagline "Pop your way to massive winnings!"'

# 3a. Build the finished italic text in a disposable scratch paragraph.
$italicIndex = $d.Paragraphs.Count
$sourceForImg = $d.Paragraphs($italicIndex).Range
$sourceForImgNoMark = $d.Range($sourceForImg.Start, $sourceForImg.End - 1)
$sourceForImgNoMark.Copy()

$scratchIndex2 = New-ScratchParagraph
$scratchRange2 = $d.Paragraphs($scratchIndex2).Range
$scratchCollapsed2 = $d.Range($scratchRange2.Start, $scratchRange2.Start)
$scratchCollapsed2.Paste()

$scratchRange2 = $d.Paragraphs($scratchIndex2).Range
[void]$scratchRange2.Find.Execute(
    "Read our review of CherryPop, an online slot game with PopWins mode and 5 free spins. Play it for free!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newImagePrompt,
    2)

# 3b. Copy the finished scratch paragraph's content (no paragraph mark).
$scratchRange2 = $d.Paragraphs($scratchIndex2).Range
$scratchNoMark2 = $d.Range($scratchRange2.Start, $scratchRange2.End - 1)
$scratchNoMark2.Copy()

# 3c. Clear the real target paragraph's existing text (collapses its runs
#     down to a single clean, un-formatted <w:r/>) ...
$italicIndex = $d.Paragraphs.Count - 1
$targetRange = $d.Paragraphs($italicIndex).Range
$targetTextLen = $targetRange.End - $targetRange.Start - 1
$targetTextRange = $d.Range($targetRange.Start, $targetRange.Start + $targetTextLen)
$targetTextRange.Delete()

# 3d. ... then paste the finished (already-italic) content into it.
$targetRange = $d.Paragraphs($italicIndex).Range
$targetCollapsed = $d.Range($targetRange.Start, $targetRange.Start)
$targetCollapsed.Paste()

# 3e. Discard the scratch paragraph again.
$scratchIndex2 = $d.Paragraphs.Count
$d.Paragraphs($scratchIndex2).Range.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
